# Peru Liga 1 workbook update (06-04-2024 01:36)
# - Re-sorts rows 180-188 (a handful of matches get re-ordered / corrected ids)
# - Corrects / fills in row 272 (date + FTHG/FTAG/FTR result that had arrived)
# - Appends three brand-new upcoming fixtures as rows 273-275

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Rows 180-188: permute the match rows (column A, the running id, stays put
#    - only B..AC travel with the record they describe).
#    after[row] = before[mapping[row]]
# ---------------------------------------------------------------------------

$snap180 = $ws.Range("B180:AC180").Value2
$snap181 = $ws.Range("B181:AC181").Value2
$snap182 = $ws.Range("B182:AC182").Value2
$snap183 = $ws.Range("B183:AC183").Value2
$snap184 = $ws.Range("B184:AC184").Value2
$snap185 = $ws.Range("B185:AC185").Value2
$snap186 = $ws.Range("B186:AC186").Value2
$snap187 = $ws.Range("B187:AC187").Value2
$snap188 = $ws.Range("B188:AC188").Value2

$ws.Range("B180:AC180").Value2 = $snap181
$ws.Range("B181:AC181").Value2 = $snap180
$ws.Range("B182:AC182").Value2 = $snap182
$ws.Range("B183:AC183").Value2 = $snap185
$ws.Range("B184:AC184").Value2 = $snap187
$ws.Range("B185:AC185").Value2 = $snap188
$ws.Range("B186:AC186").Value2 = $snap183
$ws.Range("B187:AC187").Value2 = $snap184
$ws.Range("B188:AC188").Value2 = $snap186

# ---------------------------------------------------------------------------
# 2) Row 272: this fixture has now been played - correct the kickoff time and
#    fill in the result columns (F/G swap to the correct home/away teams,
#    FTHG/FTAG/FTR appear, and the odds refresh).
# ---------------------------------------------------------------------------

$ws.Cells.Item(272, 2).Value2 = 8011505          # B272 id
$ws.Cells.Item(272, 5).Value2 = 45387.6875        # E272 Date
$ws.Cells.Item(272, 6).Value2 = "Atletico Grau"   # F272 HomeTeam
$ws.Cells.Item(272, 7).Value2 = "Sport Boys"      # G272 AwayTeam
$ws.Cells.Item(272, 8).Value2 = 0                 # H272 FTHG
$ws.Cells.Item(272, 9).Value2 = 0                 # I272 FTAG
$ws.Cells.Item(272, 10).Value2 = "D"              # J272 FTR
$ws.Cells.Item(272, 11).Value2 = 1.533            # K272 oddH_op
$ws.Cells.Item(272, 12).Value2 = 4                # L272 oddD_op
$ws.Cells.Item(272, 13).Value2 = 5.5              # M272 oddA_op
$ws.Cells.Item(272, 14).Value2 = 1.533            # N272 oddH
$ws.Cells.Item(272, 15).Value2 = 3.8              # O272 oddD
$ws.Cells.Item(272, 16).Value2 = 5.75              # P272 oddA
$ws.Cells.Item(272, 17).Value2 = -1               # Q272 Ah
$ws.Cells.Item(272, 18).Value2 = 1.925            # R272 oddAHH
$ws.Cells.Item(272, 19).Value2 = 1.875            # S272 oddAHA
$ws.Cells.Item(272, 20).Value2 = 2.5              # T272 AhOU
$ws.Cells.Item(272, 21).Value2 = 1.95             # U272 oddAHOver
$ws.Cells.Item(272, 22).Value2 = 1.85             # V272 oddAHUnder
$ws.Cells.Item(272, 23).Value2 = -1               # W272 PLH
$ws.Cells.Item(272, 24).Value2 = 2.8              # X272 PLD
$ws.Cells.Item(272, 25).Value2 = -1               # Y272 PLA
$ws.Cells.Item(272, 26).Value2 = -1               # Z272 PL_Ahh
$ws.Cells.Item(272, 27).Value2 = 0.875            # AA272 PL_Aha
$ws.Cells.Item(272, 28).Value2 = -1               # AB272 PL_AhOver
$ws.Cells.Item(272, 29).Value2 = 0.8500000000000001 # AC272 PL_AhUnder

# ---------------------------------------------------------------------------
# 3) Rows 273-275: brand-new upcoming fixtures appended at the bottom. They
#    have not been played yet so FTHG/FTAG/FTR (H/I/J) stay blank, matching
#    the original layout used for not-yet-played matches.
#    Bring over the id-column / date-column formatting from row 272 first.
# ---------------------------------------------------------------------------

$ws.Range("A272").Copy() | Out-Null
$ws.Range("A273:A275").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("E272").Copy() | Out-Null
$ws.Range("E273:E275").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$excel.CutCopyMode = 0

# --- Row 273 --------------------------------------------------------------
$ws.Cells.Item(273, 1).Value2 = 271
$ws.Cells.Item(273, 2).Value2 = 8011507
$ws.Cells.Item(273, 3).Value2 = "Peru Liga 1"
$ws.Cells.Item(273, 4).Value2 = "Peru Liga 1"
$ws.Cells.Item(273, 5).Value2 = 45388.625
$ws.Cells.Item(273, 6).Value2 = "CD Los Chankas"
$ws.Cells.Item(273, 7).Value2 = "FBC Melgar"
$ws.Cells.Item(273, 11).Value2 = 3
$ws.Cells.Item(273, 12).Value2 = 3.2
$ws.Cells.Item(273, 13).Value2 = 2.3
$ws.Cells.Item(273, 14).Value2 = 3.5
$ws.Cells.Item(273, 15).Value2 = 3.2
$ws.Cells.Item(273, 16).Value2 = 2.05
$ws.Cells.Item(273, 17).Value2 = 0.25
$ws.Cells.Item(273, 18).Value2 = 2.05
$ws.Cells.Item(273, 19).Value2 = 1.8
$ws.Cells.Item(273, 20).Value2 = 2.5
$ws.Cells.Item(273, 21).Value2 = 1.925
$ws.Cells.Item(273, 22).Value2 = 1.925
$ws.Cells.Item(273, 23).Value2 = 0
$ws.Cells.Item(273, 24).Value2 = 0
$ws.Cells.Item(273, 25).Value2 = 0
$ws.Cells.Item(273, 26).Value2 = 0
$ws.Cells.Item(273, 27).Value2 = 0

# --- Row 274 --------------------------------------------------------------
$ws.Cells.Item(274, 1).Value2 = 272
$ws.Cells.Item(274, 2).Value2 = 8011508
$ws.Cells.Item(274, 3).Value2 = "Peru Liga 1"
$ws.Cells.Item(274, 4).Value2 = "Peru Liga 1"
$ws.Cells.Item(274, 5).Value2 = 45388.72916666666
$ws.Cells.Item(274, 6).Value2 = "Union Comercio"
$ws.Cells.Item(274, 7).Value2 = "Cesar Vallejo"
$ws.Cells.Item(274, 11).Value2 = 2.2
$ws.Cells.Item(274, 12).Value2 = 3.3
$ws.Cells.Item(274, 13).Value2 = 3.1
$ws.Cells.Item(274, 14).Value2 = 2.55
$ws.Cells.Item(274, 15).Value2 = 3.3
$ws.Cells.Item(274, 16).Value2 = 2.4
$ws.Cells.Item(274, 17).Value2 = 0
$ws.Cells.Item(274, 18).Value2 = 2
$ws.Cells.Item(274, 19).Value2 = 1.85
$ws.Cells.Item(274, 20).Value2 = 2.5
$ws.Cells.Item(274, 21).Value2 = 2.05
$ws.Cells.Item(274, 22).Value2 = 1.8
$ws.Cells.Item(274, 23).Value2 = 0
$ws.Cells.Item(274, 24).Value2 = 0
$ws.Cells.Item(274, 25).Value2 = 0
$ws.Cells.Item(274, 26).Value2 = 0
$ws.Cells.Item(274, 27).Value2 = 0

# --- Row 275 --------------------------------------------------------------
$ws.Cells.Item(275, 1).Value2 = 273
$ws.Cells.Item(275, 2).Value2 = 8012287
$ws.Cells.Item(275, 3).Value2 = "Peru Liga 1"
$ws.Cells.Item(275, 4).Value2 = "Peru Liga 1"
$ws.Cells.Item(275, 5).Value2 = 45388.8125
$ws.Cells.Item(275, 6).Value2 = "Carlos Manucci"
$ws.Cells.Item(275, 7).Value2 = "Alianza Lima"
$ws.Cells.Item(275, 11).Value2 = 3.5
$ws.Cells.Item(275, 12).Value2 = 3.4
$ws.Cells.Item(275, 13).Value2 = 2
$ws.Cells.Item(275, 14).Value2 = 4
$ws.Cells.Item(275, 15).Value2 = 3.4
$ws.Cells.Item(275, 16).Value2 = 1.75
$ws.Cells.Item(275, 17).Value2 = 0.5
$ws.Cells.Item(275, 18).Value2 = 2.05
$ws.Cells.Item(275, 19).Value2 = 1.8
$ws.Cells.Item(275, 20).Value2 = 2.5
$ws.Cells.Item(275, 21).Value2 = 1.875
$ws.Cells.Item(275, 22).Value2 = 1.975
$ws.Cells.Item(275, 23).Value2 = 0
$ws.Cells.Item(275, 24).Value2 = 0
$ws.Cells.Item(275, 25).Value2 = 0
$ws.Cells.Item(275, 26).Value2 = 0
$ws.Cells.Item(275, 27).Value2 = 0

Write-Output "Peru Liga 1 sheet updated: rows 180-188 reordered, row 272 completed, rows 273-275 appended."
